# Generate Report for Handoff
#
# Swap which row represents which source file on every sheet (the
# a9ddeb99... file moves to row 2, 0d27cd94... moves to row 3) and mark
# the 0d27cd94... file as handed off ("Ready for handoff") with fresh
# handoff timestamps / handoff-file links.

$wb = $excel.ActiveWorkbook

$mdUrl_0d27 = "https://github.com/OpenLocalizationTest/oltest/blob/26716da358451b61a1a12fa209bf516df58c3927/e2e/0d27cd94-d15d-4bf9-b30a-0e4587535603.md"
$mdUrl_a9dd = "https://github.com/OpenLocalizationTest/oltest/blob/26716da358451b61a1a12fa209bf516df58c3927/e2e/a9ddeb99-0341-40e1-97b3-424324c8b84d.md"

$zhUrl_0d27 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2ae1007edcd820295040597027568884690d100/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0d27cd94-d15d-4bf9-b30a-0e4587535603.0c4d2487f5ba1fb09dc064a42b1223a44bcb5e7c.zh-cn.xlf"
$zhUrl_a9dd = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2ae1007edcd820295040597027568884690d100/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a9ddeb99-0341-40e1-97b3-424324c8b84d.aadca84f17d073cd11d4a7975d923e8ea0e98318.zh-cn.xlf"

$deUrl_0d27 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a53fdc11f1d3ebf9b763c075b8107a7e7d28576/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0d27cd94-d15d-4bf9-b30a-0e4587535603.0c4d2487f5ba1fb09dc064a42b1223a44bcb5e7c.de-de.xlf"
$deUrl_a9dd = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a53fdc11f1d3ebf9b763c075b8107a7e7d28576/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a9ddeb99-0341-40e1-97b3-424324c8b84d.aadca84f17d073cd11d4a7975d923e8ea0e98318.de-de.xlf"

$name_0d27 = "0d27cd94-d15d-4bf9-b30a-0e4587535603.md"
$name_a9dd = "a9ddeb99-0341-40e1-97b3-424324c8b84d.md"

$miss = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "2016-15-19 06:15:54"

$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-16-19 06:16:32"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl_0d27, $miss, $miss, $name_a9dd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl_a9dd, $miss, $miss, $name_0d27) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C2").Value = "In Translation"
$ws.Range("E2").Value = "2016-03-19 06:15:51"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-19 06:16:29"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl_0d27, $miss, $miss, $name_a9dd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl_0d27, $miss, $miss, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $zhUrl_0d27, $miss, $miss, "a9ddeb99-0341-40e1-97b3-424324c8b84d.aadca84f17d073cd11d4a7975d923e8ea0e98318.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl_a9dd, $miss, $miss, $name_0d27) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl_a9dd, $miss, $miss, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $zhUrl_a9dd, $miss, $miss, "0d27cd94-d15d-4bf9-b30a-0e4587535603.0c4d2487f5ba1fb09dc064a42b1223a44bcb5e7c.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C2").Value = "In Translation"
$ws.Range("E2").Value = "2016-03-19 06:15:54"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-19 06:16:32"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl_0d27, $miss, $miss, $name_a9dd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl_0d27, $miss, $miss, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $deUrl_0d27, $miss, $miss, "a9ddeb99-0341-40e1-97b3-424324c8b84d.aadca84f17d073cd11d4a7975d923e8ea0e98318.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl_a9dd, $miss, $miss, $name_0d27) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl_a9dd, $miss, $miss, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $deUrl_a9dd, $miss, $miss, "0d27cd94-d15d-4bf9-b30a-0e4587535603.0c4d2487f5ba1fb09dc064a42b1223a44bcb5e7c.de-de.xlf") | Out-Null

$wb.Save()
